# Update the scrape timestamp for every data row (2-130)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2:O130").Value = "2022-07-29 20:58:31"

# Rows 31-33 got re-ordered in the new scrape (same three products, new row order):
#   new row 31 <- old row 32 (Avela ... Noir ... 6075749004)
#   new row 32 <- old row 33 (Naturaline Herren T-Shirt ... 6031467007)
#   new row 33 <- old row 31 (Avela ... Hasel ... 6075749002)
# Columns A (id) and H (price) hold numbers-as-text in this sheet, so force
# text formatting before writing them or Excel will coerce them to numeric.
# (H31 keeps its original "1.95" value, so it is left untouched.)

$ws.Range("A31:A33").NumberFormat = "@"
$ws.Range("H32:H33").NumberFormat = "@"

$ws.Range("A31").Value = "6075749004"
$ws.Range("B31").Value = "Avela Strumpfhose Top Size Noir  13 - 14"
$ws.Range("C31").Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-top-size-noir-13-14/p/6075749004"
$ws.Range("M31").Value = "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'struempfe']"
$ws.Range("N31").Value = "Avela Strumpfhose Top Size Noir  13 - 14 50% Aktion 1.95 Schweizer Franken statt 3.95 Schweizer Franken"

$ws.Range("A32").Value = "6031467007"
$ws.Range("B32").Value = "Naturaline Herren T-Shirt Kurzarm weiss M"
$ws.Range("C32").Value = "/de/haushalt-tier/bekleidung/shirts-pullover/herren-shirt/naturaline-herren-t-shirt-kurzarm-weiss-m/p/6031467007"
$ws.Range("G32").Value = "Coop"
$ws.Range("H32").Value = "24.95"
$ws.Range("M32").Value = "['haushalt-tier', 'bekleidung', 'shirts-pullover', 'herren-shirt']"
$ws.Range("N32").Value = "Naturaline Herren T-Shirt Kurzarm weiss M 24.95 Schweizer Franken"

$ws.Range("A33").Value = "6075749002"
$ws.Range("B33").Value = "Avela Strumpfhose Top Size Hasel  13 - 14"
$ws.Range("C33").Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-top-size-hasel-13-14/p/6075749002"
$ws.Range("G33").Value = "Avela"
$ws.Range("H33").Value = "1.95"
$ws.Range("M33").Value = "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'struempfe']"
$ws.Range("N33").Value = "Avela Strumpfhose Top Size Hasel  13 - 14 50% Aktion 1.95 Schweizer Franken statt 3.95 Schweizer Franken"
